$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.709516000000001
$ws.Range("H2").Value = 17.128548
$ws.Range("I2").Value = 0.394755842864453
$ws.Range("J2").Value = 0.394755842864453
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 39.413899
$ws.Range("N2").Value = 118.241697
$ws.Range("O2").Value = 0.7310109930444597
$ws.Range("P2").Value = 0.7310109930444597
$ws.Range("Q2").Value = 225.034286962884
$ws.Range("R2").Value = 2025.308582665956
$ws.Range("S2").Value = 0.2885708607024465
$ws.Range("T2").Value = 0.2885708607024465

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.709516000000001
$ws.Range("H3").Value = 17.128548
$ws.Range("I3").Value = 0.394755842864453
$ws.Range("J3").Value = 0.394755842864453
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 11.40791933333333
$ws.Range("N3").Value = 34.223758
$ws.Range("O3").Value = 0.2115830874897988
$ws.Range("P3").Value = 0.2115830874897987
$ws.Range("Q3").Value = 65.13369796037601
$ws.Range("R3").Value = 586.2032816433841
$ws.Range("S3").Value = 0.08352366003789881
$ws.Range("T3").Value = 0.0835236600378988

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.709516000000001
$ws.Range("H4").Value = 17.128548
$ws.Range("I4").Value = 0.394755842864453
$ws.Range("J4").Value = 0.394755842864453
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.095153333333334
$ws.Range("N4").Value = 9.28546
$ws.Range("O4").Value = 0.05740591946574151
$ws.Range("P4").Value = 0.0574059194657415
$ws.Range("Q4").Value = 17.67182747912
$ws.Range("R4").Value = 159.04644731208
$ws.Range("S4").Value = 0.0226613221241077
$ws.Range("T4").Value = 0.02266132212410769

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.281282666666666
$ws.Range("H5").Value = 9.843848
$ws.Range("I5").Value = 0.2268678299100168
$ws.Range("J5").Value = 0.2268678299100168
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 39.413899
$ws.Range("N5").Value = 118.241697
$ws.Range("O5").Value = 0.7310109930444597
$ws.Range("P5").Value = 0.7310109930444597
$ws.Range("Q5").Value = 129.3281436144507
$ws.Range("R5").Value = 1163.953292530056
$ws.Range("S5").Value = 0.165842877632363
$ws.Range("T5").Value = 0.165842877632363

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.281282666666666
$ws.Range("H6").Value = 9.843848
$ws.Range("I6").Value = 0.2268678299100168
$ws.Range("J6").Value = 0.2268678299100168
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 11.40791933333333
$ws.Range("N6").Value = 34.223758
$ws.Range("O6").Value = 0.2115830874897988
$ws.Range("P6").Value = 0.2115830874897987
$ws.Range("Q6").Value = 37.43260797119822
$ws.Range("R6").Value = 336.893471740784
$ws.Range("S6").Value = 0.04800139590447187
$ws.Range("T6").Value = 0.04800139590447187

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.281282666666666
$ws.Range("H7").Value = 9.843848
$ws.Range("I7").Value = 0.2268678299100168
$ws.Range("J7").Value = 0.2268678299100168
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.095153333333334
$ws.Range("N7").Value = 9.28546
$ws.Range("O7").Value = 0.05740591946574151
$ws.Range("P7").Value = 0.0574059194657415
$ws.Range("Q7").Value = 10.15607298334222
$ws.Range("R7").Value = 91.40465685008
$ws.Range("S7").Value = 0.01302355637318197
$ws.Range("T7").Value = 0.01302355637318197

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.388848333333333
$ws.Range("H8").Value = 4.166545
$ws.Range("I8").Value = 0.09602495105292475
$ws.Range("J8").Value = 0.09602495105292475
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 39.413899
$ws.Range("N8").Value = 118.241697
$ws.Range("O8").Value = 0.7310109930444597
$ws.Range("P8").Value = 0.7310109930444597
$ws.Range("Q8").Value = 54.73992793631834
$ws.Range("R8").Value = 492.659351426865
$ws.Range("S8").Value = 0.07019529482624416
$ws.Range("T8").Value = 0.07019529482624416

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.388848333333333
$ws.Range("H9").Value = 4.166545
$ws.Range("I9").Value = 0.09602495105292475
$ws.Range("J9").Value = 0.09602495105292475
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.40791933333333
$ws.Range("N9").Value = 34.223758
$ws.Range("O9").Value = 0.2115830874897988
$ws.Range("P9").Value = 0.2115830874897987
$ws.Range("Q9").Value = 15.84386975290111
$ws.Range("R9").Value = 142.59482777611
$ws.Range("S9").Value = 0.02031725561983462
$ws.Range("T9").Value = 0.02031725561983462

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.388848333333333
$ws.Range("H10").Value = 4.166545
$ws.Range("I10").Value = 0.09602495105292475
$ws.Range("J10").Value = 0.09602495105292475
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.095153333333334
$ws.Range("N10").Value = 9.28546
$ws.Range("O10").Value = 0.05740591946574151
$ws.Range("P10").Value = 0.0574059194657415
$ws.Range("Q10").Value = 4.298698548411112
$ws.Range("R10").Value = 38.68828693570001
$ws.Range("S10").Value = 0.005512400606845969
$ws.Range("T10").Value = 0.005512400606845967

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.083764
$ws.Range("H11").Value = 12.251292
$ws.Range("I11").Value = 0.2823513761726055
$ws.Range("J11").Value = 0.2823513761726055
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 39.413899
$ws.Range("N11").Value = 118.241697
$ws.Range("O11").Value = 0.7310109930444597
$ws.Range("P11").Value = 0.7310109930444597
$ws.Range("Q11").Value = 160.957061835836
$ws.Range("R11").Value = 1448.613556522524
$ws.Range("S11").Value = 0.2064019598834061
$ws.Range("T11").Value = 0.2064019598834061

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.083764
$ws.Range("H12").Value = 12.251292
$ws.Range("I12").Value = 0.2823513761726055
$ws.Range("J12").Value = 0.2823513761726055
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 11.40791933333333
$ws.Range("N12").Value = 34.223758
$ws.Range("O12").Value = 0.2115830874897988
$ws.Range("P12").Value = 0.2115830874897987
$ws.Range("Q12").Value = 46.58725028837068
$ws.Range("R12").Value = 419.2852525953361
$ws.Range("S12").Value = 0.05974077592759346
$ws.Range("T12").Value = 0.05974077592759346

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.083764
$ws.Range("H13").Value = 12.251292
$ws.Range("I13").Value = 0.2823513761726055
$ws.Range("J13").Value = 0.2823513761726055
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.095153333333334
$ws.Range("N13").Value = 9.28546
$ws.Range("O13").Value = 0.05740591946574151
$ws.Range("P13").Value = 0.0574059194657415
$ws.Range("Q13").Value = 12.63987575714667
$ws.Range("R13").Value = 113.75888181432
$ws.Range("S13").Value = 0.01620864036160588
$ws.Range("T13").Value = 0.01620864036160587
